# REFACTOR: Removing unwanted suites
#
# "NextAndBackButtons" (row 4) is an unwanted/obsolete test suite being
# removed from the Step1 statistics sheet. Its row is dropped, the rows
# below it shift up to fill the gap, and every range/formula that used to
# span down to the old last data row (36/37/50) is tightened by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Row 4 becomes what used to be row 5 (ProgressBar / Suited to Manual) ---
$ws.Range("A4").Value = "ProgressBar"
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = "Suited to Manual"

# --- Row 5 becomes what used to be row 6 (RequiredFields / Automated) ---
$ws.Range("A5").Value = "RequiredFields"
$ws.Range("B5").Value = 4
$ws.Range("D5").Value = "Automated"

# --- Old row 6 entries (now duplicated into row 5) are gone entirely ---
$ws.Range("A6:D6").Clear()

# --- Aggregate formulas shrink their ranges by one row ---
$ws.Range("E1").Formula = '=COUNTA($A$2:A36)'
$ws.Range("G1").Formula = '=COUNTIF($D$2:D35,"Ready to Write")+COUNTIF($D$2:D35,"Outdated")'
$ws.Range("G2").Formula = '=COUNTIF($D$2:D35,"Automated")+COUNTIF($D$2:D35,"Finished")'
$ws.Range("G4").Formula = '=SUM($C$2:C35)'
$ws.Range("G5").Formula = '=SUM($B$2:B35)'

# --- "Blank check" conditional format range tightens the same way ---
$fcs = $ws.Range("D2:D50").FormatConditions
$blankRule = $fcs.Item(1)
$blankRule.ModifyAppliesToRange($ws.Range("D2:D49"))

# --- Re-balance the "containsText" rule priorities/captions ---
$fTesting = $fcs.Item(2)
$fTesting.Text = "Finished"
$fTesting.Formula1 = '=NOT(ISERROR(SEARCH("Finished",D1)))'
$fTesting.Priority = 1

$fWriting = $fcs.Item(3)
$fWriting.Text = "Automated"
$fWriting.Formula1 = '=NOT(ISERROR(SEARCH("Automated",D1)))'
$fWriting.Priority = 2

$fAutomated = $fcs.Item(5)
$fAutomated.Text = "Testing"
$fAutomated.Formula1 = '=NOT(ISERROR(SEARCH("Testing",D1)))'
$fAutomated.Priority = 4

$fFinished = $fcs.Item(6)
$fFinished.Text = "Writing"
$fFinished.Formula1 = '=NOT(ISERROR(SEARCH("Writing",D1)))'
$fFinished.Priority = 5

# --- Leave the cursor where the editor left it ---
[void]$ws.Range("D12").Select()
